$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 119.25
$ws.Range("I9").Value = 119.25
$ws.Range("K9").Value = 119.25
$ws.Range("M9").Value = 49.75

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 529.2778
$ws.Range("I33").Value = 575.1667
$ws.Range("J33").Value = 437.5
$ws.Range("K33").Value = 575.1667
$ws.Range("L33").Value = 437.5
$ws.Range("M33").Value = -346.1667
$ws.Range("N33").Value = -895.5

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2400.6
$ws.Range("I40").Value = 2250.75
$ws.Range("K40").Value = 2250.75
$ws.Range("M40").Value = -2075.75

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5363.4546
$ws.Range("I76").Value = 5800
$ws.Range("K76").Value = 5800
$ws.Range("M76").Value = -5485

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 5363.4546
$ws.Range("I79").Value = 5800
$ws.Range("K79").Value = 5800
$ws.Range("M79").Value = -4708

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1216.9231
$ws.Range("I100").Value = 760
$ws.Range("J100").Value = 2740
$ws.Range("K100").Value = 760
$ws.Range("L100").Value = 2740
$ws.Range("M100").Value = -219
$ws.Range("N100").Value = -3822

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 1254

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3612.6155
$ws.Range("I116").Value = 2991.2
$ws.Range("J116").Value = 4001
$ws.Range("K116").Value = 2991.2
$ws.Range("L116").Value = 4001
$ws.Range("M116").Value = 450.8000000000002
$ws.Range("N116").Value = -10885

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7414853
$ws.Range("I132").Value = 9013910
$ws.Range("K132").Value = 27041730
$ws.Range("M132").Value = -27039200

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1439.7869
$ws.Range("I137").Value = 1365.2941
$ws.Range("K137").Value = 4095.8823
$ws.Range("M137").Value = -1545.8823

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 6675.4546
$ws.Range("I28").Value = 6675.4546
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 6675.4546
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -6483.4546
$ws.Range("N28").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5433.636
$ws.Range("I31").Value = 5433.636
$ws.Range("K31").Value = 5433.636
$ws.Range("M31").Value = -5139.636

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3988.8
$ws.Range("I32").Value = 3378.161
$ws.Range("J32").Value = 8075.385
$ws.Range("K32").Value = 3378.161
$ws.Range("L32").Value = 8075.385
$ws.Range("M32").Value = -3091.161
$ws.Range("N32").Value = -8649.385

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1410.8667
$ws.Range("I45").Value = 1270.375
$ws.Range("J45").Value = 1571.4286
$ws.Range("K45").Value = 1270.375
$ws.Range("L45").Value = 1571.4286
$ws.Range("M45").Value = -893.375
$ws.Range("N45").Value = -2325.4286

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 71429850
$ws.Range("I61").Value = 90910060
$ws.Range("K61").Value = 90910060
$ws.Range("M61").Value = -90909848

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2507.3076
$ws.Range("I74").Value = 2113.1428
$ws.Range("J74").Value = 4162.8
$ws.Range("K74").Value = 2113.1428
$ws.Range("L74").Value = 4162.8
$ws.Range("M74").Value = -1239.1428
$ws.Range("N74").Value = -5910.8

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2507.3076
$ws.Range("I77").Value = 2113.1428
$ws.Range("J77").Value = 4162.8
$ws.Range("K77").Value = 10565.714
$ws.Range("L77").Value = 20814
$ws.Range("M77").Value = -6197.714
$ws.Range("N77").Value = -29550

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 6675.4546
$ws.Range("I99").Value = 6675.4546
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6675.4546
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3680.4546
$ws.Range("N99").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 71429850
$ws.Range("I136").Value = 90910060
$ws.Range("K136").Value = 272730180
$ws.Range("M136").Value = -272727630

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 20327.75
$ws.Range("I97").Value = 5437
$ws.Range("K97").Value = 5437
$ws.Range("M97").Value = -4446

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 111112240
$ws.Range("I99").Value = 125001070
$ws.Range("J99").Value = 1600
$ws.Range("K99").Value = 125001070
$ws.Range("L99").Value = 1600
$ws.Range("M99").Value = -124999572
$ws.Range("N99").Value = -4596

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1150.7222
$ws.Range("I107").Value = 1142.8572
$ws.Range("K107").Value = 1142.8572
$ws.Range("M107").Value = 777.1428000000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1518.8572
$ws.Range("I134").Value = 1355.4166
$ws.Range("K134").Value = 4066.2498
$ws.Range("M134").Value = -1531.2498

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1609.6
$ws.Range("I99").Value = 1524.4
$ws.Range("J99").Value = 1780
$ws.Range("K99").Value = 1524.4
$ws.Range("L99").Value = 1780
$ws.Range("M99").Value = -26.40000000000009
$ws.Range("N99").Value = -4776

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 895.2222
$ws.Range("I107").Value = 470.30768
$ws.Range("K107").Value = 470.30768
$ws.Range("M107").Value = 1449.69232

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1609.6
$ws.Range("I126").Value = 1524.4
$ws.Range("J126").Value = 1780
$ws.Range("K126").Value = 4573.200000000001
$ws.Range("L126").Value = 5340
$ws.Range("M126").Value = -2103.200000000001
$ws.Range("N126").Value = -10280

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4139.4546
$ws.Range("J39").Value = 4214.8887
$ws.Range("L39").Value = 12644.6661
$ws.Range("N39").Value = -13232.6661

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 469.9375
$ws.Range("I114").Value = 331.2857
$ws.Range("J114").Value = 577.7778
$ws.Range("K114").Value = 993.8571000000001
$ws.Range("L114").Value = 1733.3334
$ws.Range("M114").Value = 2260.1429
$ws.Range("N114").Value = -8241.3334

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 5688
$ws.Range("I120").Value = 4999.5
$ws.Range("K120").Value = 14998.5
$ws.Range("M120").Value = -10160.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 26788482
$ws.Range("I137").Value = 83335030
$ws.Range("J137").Value = 3275.3684
$ws.Range("K137").Value = 250005090
$ws.Range("L137").Value = 9826.1052
$ws.Range("M137").Value = -249999990
$ws.Range("N137").Value = -20026.1052

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 2407.3208
$ws.Range("I138").Value = 2633.7693
$ws.Range("J138").Value = 2333.725
$ws.Range("K138").Value = 7901.3079
$ws.Range("L138").Value = 7001.174999999999
$ws.Range("M138").Value = -2761.3079
$ws.Range("N138").Value = -17281.175

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50003480
$ws.Range("I70").Value = 50003690
$ws.Range("J70").Value = 50003220
$ws.Range("K70").Value = 50003690
$ws.Range("L70").Value = 50003220
$ws.Range("M70").Value = -50003420
$ws.Range("N70").Value = -50003760

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 50003480
$ws.Range("I73").Value = 50003690
$ws.Range("J73").Value = 50003220
$ws.Range("K73").Value = 50003690
$ws.Range("L73").Value = 50003220
$ws.Range("M73").Value = -50002754
$ws.Range("N73").Value = -50005092

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 10471
$ws.Range("I99").Value = 10471
$ws.Range("K99").Value = 10471
$ws.Range("M99").Value = -8225

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 7166.5
$ws.Range("J109").Value = 7166.5
$ws.Range("L109").Value = 7166.5
$ws.Range("N109").Value = -9246.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 12581.044
$ws.Range("J136").Value = 12581.044
$ws.Range("L136").Value = 37743.132
$ws.Range("N136").Value = -42843.132

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1271.6364
$ws.Range("I68").Value = 1271.6364
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1271.6364
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -522.6364000000001
$ws.Range("N68").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1271.6364
$ws.Range("I71").Value = 1271.6364
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 6358.182000000001
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -2614.182000000001
$ws.Range("N71").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1768.6364
$ws.Range("I82").Value = 1753.8572
$ws.Range("J82").Value = 1794.5
$ws.Range("K82").Value = 1753.8572
$ws.Range("L82").Value = 1794.5
$ws.Range("M82").Value = -1392.8572
$ws.Range("N82").Value = -2516.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1768.6364
$ws.Range("I85").Value = 1753.8572
$ws.Range("J85").Value = 1794.5
$ws.Range("K85").Value = 1753.8572
$ws.Range("L85").Value = 1794.5
$ws.Range("M85").Value = -505.8571999999999
$ws.Range("N85").Value = -4290.5

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 6000
$ws.Range("J9").Value = 6000
$ws.Range("L9").Value = 6000
$ws.Range("N9").Value = -6280

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9000
$ws.Range("J69").Value = 9000
$ws.Range("L69").Value = 9000
$ws.Range("N69").Value = -10498

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 9000
$ws.Range("J72").Value = 9000
$ws.Range("L72").Value = 27000
$ws.Range("N72").Value = -34488

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1408.6
$ws.Range("I81").Value = 1670.5
$ws.Range("J81").Value = 1234
$ws.Range("K81").Value = 3341
$ws.Range("L81").Value = 2468
$ws.Range("M81").Value = -2280
$ws.Range("N81").Value = -4590

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1408.6
$ws.Range("I84").Value = 1670.5
$ws.Range("J84").Value = 1234
$ws.Range("K84").Value = 16705
$ws.Range("L84").Value = 12340
$ws.Range("M84").Value = -11401
$ws.Range("N84").Value = -22948

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4549.9287
$ws.Range("I96").Value = 4069.9
$ws.Range("J96").Value = 5750
$ws.Range("K96").Value = 4069.9
$ws.Range("L96").Value = 5750
$ws.Range("M96").Value = -2696.9
$ws.Range("N96").Value = -8496

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 58824704
$ws.Range("I126").Value = 100000500
$ws.Range("J126").Value = 2146.2856
$ws.Range("K126").Value = 300001500
$ws.Range("L126").Value = 6438.8568
$ws.Range("M126").Value = -299999030
$ws.Range("N126").Value = -11378.8568

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2170.861
$ws.Range("I132").Value = 1989.1562
$ws.Range("K132").Value = 5967.4686
$ws.Range("M132").Value = -3437.4686
